$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting all existing data down by one row.
$ws.Rows.Item(1).Insert()

# Write the new header row.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Descrição"
